$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.693.30'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '3.728.32'
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.09'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.49'
$ws.Range("E6").Value = '  -4.46%  '
$ws.Range("D7").Value = '3.726.37'
$ws.Range("E7").Value = '  -1.14%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("E10").Value = '  +2.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.39'
$ws.Range("E11").Value = '  +3.13%  '
$ws.Range("E12").Value = '  -1.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.95'
$ws.Range("E13").Value = '  -1.40%  '
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '4.355.02'
$ws.Range("E15").Value = '  -1.13%  '
$ws.Range("D16").Value = '3.730.13'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").Value = '68.687.91'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("E20").Value = '  +4.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.97'
$ws.Range("E21").Value = '  +0.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.04'
$ws.Range("E22").Value = '  +10.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.721'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.84'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("E25").Value = '  -5.13%  '
$ws.Range("E26").Value = '  -2.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.36'
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("E28").Value = '  -1.19%  '
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  -0.61%  '
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.60'
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").Value = '3.876.80'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").Value = '3.663.63'
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("E36").Value = '  -1.20%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  +0.04%  '
$ws.Range("E40").Value = '  -1.75%  '
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.73'
$ws.Range("E42").Value = '  -3.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.13'
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("E44").Value = '  -1.95%  '
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  +0.65%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.66'
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.88'
$ws.Range("E49").Value = '  +1.41%  '
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").Value = '2.739.55'
$ws.Range("E51").Value = '  -3.35%  '
